$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: classical-best-embeddings -> classical-best-embed (in A2 only)
$ws.Range("A2").Value = "classical-best-embed vs. classical-best-tfidf"
$ws.Range("C2").Value = 0.079
$ws.Range("D2").Value = 0.062
$ws.Range("I2").Value = 0.061
$ws.Range("J2").Value = 0.065

# Row 3: text unchanged, values updated
$ws.Range("C3").Value = 0.089
$ws.Range("D3").Value = 0.138
$ws.Range("E3").Value = 0.129
$ws.Range("F3").Value = 0.119
$ws.Range("G3").Value = 0.16
$ws.Range("H3").Value = 0.187

# Row 4: BERT-base vs. classical-best-embeddings -> BERT-base vs. classical-best-embed
$ws.Range("A4").Value = "BERT-base vs. classical-best-embed"
$ws.Range("C4").Value = 0.01
$ws.Range("D4").Value = 0.076
$ws.Range("E4").Value = 0.075
$ws.Range("F4").Value = 0.07099999999999999
$ws.Range("G4").Value = 0.101
$ws.Range("H4").Value = 0.096
$ws.Range("I4").Value = 0.058
$ws.Range("J4").Value = 0.07199999999999999

# Row 5: text unchanged, values updated
$ws.Range("B5").Value = 0.437
$ws.Range("C5").Value = 0.215
$ws.Range("D5").Value = 0.194
$ws.Range("E5").Value = 0.178
$ws.Range("F5").Value = 0.16
$ws.Range("G5").Value = 0.177
$ws.Range("H5").Value = 0.196
$ws.Range("I5").Value = 0.187
$ws.Range("J5").Value = 0.187

# Row 6: BERT-base-nli vs. classical-best-embeddings -> BERT-base-nli vs. classical-best-embed
$ws.Range("A6").Value = "BERT-base-nli vs. classical-best-embed"
$ws.Range("B6").Value = 0.437
$ws.Range("C6").Value = 0.136
$ws.Range("D6").Value = 0.132
$ws.Range("E6").Value = 0.124
$ws.Range("F6").Value = 0.112
$ws.Range("G6").Value = 0.118
$ws.Range("H6").Value = 0.105
$ws.Range("I6").Value = 0.126
$ws.Range("J6").Value = 0.121

# Row 7: text unchanged, values updated
$ws.Range("B7").Value = 0.437
$ws.Range("C7").Value = 0.126
$ws.Range("D7").Value = 0.056
$ws.Range("E7").Value = 0.049
$ws.Range("F7").Value = 0.041
$ws.Range("G7").Value = 0.017
$ws.Range("H7").Value = 0.008999999999999999
$ws.Range("I7").Value = 0.068
$ws.Range("J7").Value = 0.05
